$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prepend ss_ids DOI-status annotations to comment column (H) for rows with multiple ss_ids
$ws.Range("H3").Value2 = '[ss_ids: same DOI] ' + $ws.Range("H3").Value2
$ws.Range("H5").Value2 = '[ss_ids: same DOI] ' + $ws.Range("H5").Value2
$ws.Range("H7").Value2 = '[ss_ids: no DOIs!] ' + $ws.Range("H7").Value2
$ws.Range("H8").Value2 = '[ss_ids: same DOI] ' + $ws.Range("H8").Value2
$ws.Range("H9").Value2 = '[ss_ids: multiple DOIs!] ' + $ws.Range("H9").Value2
$ws.Range("H10").Value2 = '[ss_ids: multiple DOIs!] ' + $ws.Range("H10").Value2
$ws.Range("H12").Value2 = '[ss_ids: multiple DOIs!] ' + $ws.Range("H12").Value2
$ws.Range("H13").Value2 = '[ss_ids: same DOI] ' + $ws.Range("H13").Value2
$ws.Range("H15").Value2 = '[ss_ids: same DOI] ' + $ws.Range("H15").Value2
$ws.Range("H16").Value2 = '[ss_ids: same DOI] ' + $ws.Range("H16").Value2
$ws.Range("H17").Value2 = '[ss_ids: multiple DOIs!] ' + $ws.Range("H17").Value2
$ws.Range("H19").Value2 = '[ss_ids: 1 no DOI, rest same DOI] ' + $ws.Range("H19").Value2
$ws.Range("H21").Value2 = '[ss_ids: multiple DOIs!] ' + $ws.Range("H21").Value2
$ws.Range("H24").Value2 = '[ss_ids: same DOI] ' + $ws.Range("H24").Value2
$ws.Range("H27").Value2 = '[ss_ids: same DOI] ' + $ws.Range("H27").Value2
$ws.Range("H28").Value2 = '[ss_ids: multiple DOIs!] ' + $ws.Range("H28").Value2
$ws.Range("H29").Value2 = '[ss_ids: multiple DOIs!] ' + $ws.Range("H29").Value2
$ws.Range("H33").Value2 = '[ss_ids: same DOI] ' + $ws.Range("H33").Value2
$ws.Range("H36").Value2 = '[ss_ids: same DOI] ' + $ws.Range("H36").Value2
$ws.Range("H37").Value2 = '[ss_ids: same DOI] ' + $ws.Range("H37").Value2
$ws.Range("H39").Value2 = '[ss_ids: same DOI] ' + $ws.Range("H39").Value2
$ws.Range("H44").Value2 = '[ss_ids: 1 no DOI, rest same DOI] ' + $ws.Range("H44").Value2
$ws.Range("H45").Value2 = '[ss_ids: 1 no DOI, rest same DOI] ' + $ws.Range("H45").Value2
$ws.Range("H46").Value2 = '[ss_ids: same DOI] ' + $ws.Range("H46").Value2
$ws.Range("H53").Value2 = '[ss_ids: multiple DOIs!] ' + $ws.Range("H53").Value2
$ws.Range("H55").Value2 = '[ss_ids: multiple DOIs!] ' + $ws.Range("H55").Value2
$ws.Range("H58").Value2 = '[ss_ids: same DOI] ' + $ws.Range("H58").Value2
$ws.Range("H64").Value2 = '[ss_ids: no DOIs!] ' + $ws.Range("H64").Value2
$ws.Range("H65").Value2 = '[ss_ids: same DOI] ' + $ws.Range("H65").Value2
$ws.Range("H68").Value2 = '[ss_ids: same DOI] ' + $ws.Range("H68").Value2
$ws.Range("H70").Value2 = '[ss_ids: no DOIs!] ' + $ws.Range("H70").Value2
$ws.Range("H71").Value2 = '[ss_ids: 2 no DOI, rest same DOI] ' + $ws.Range("H71").Value2
$ws.Range("H78").Value2 = '[ss_ids: multiple DOIs!] ' + $ws.Range("H78").Value2
$ws.Range("H95").Value2 = '[ss_ids: 1 no DOI, rest same DOI] ' + $ws.Range("H95").Value2
$ws.Range("H96").Value2 = '[ss_ids: 1 no DOI, rest same DOI] ' + $ws.Range("H96").Value2
$ws.Range("H99").Value2 = '[ss_ids: same DOI] ' + $ws.Range("H99").Value2
$ws.Range("H100").Value2 = '[ss_ids: same DOI] ' + $ws.Range("H100").Value2
$ws.Range("H102").Value2 = '[ss_ids: multiple DOIs!] ' + $ws.Range("H102").Value2
$ws.Range("H105").Value2 = '[ss_ids: same DOI] ' + $ws.Range("H105").Value2
$ws.Range("H108").Value2 = '[ss_ids: same DOI] ' + $ws.Range("H108").Value2

# Append new bibliography rows 121-123
$ws.Range("A120").Copy() | Out-Null
$ws.Range("A121:A123").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A121").Value2 = 1745
$ws.Range("B121").Value2 = 'Hend23a'
$ws.Range("C121").Value2 = 'article'
$ws.Range("D121").Value2 = 'Trends in the incidence of pulmonary nodules in chest computed tomography: 10-year results from two Dutch hospitals'
$ws.Range("E121").Value2 = 'Hendrix, Ward and Rutten, Matthieu and Hendrix, Nils and van Ginneken, Bram and Schaefer-Prokop, Cornelia and Scholten, Ernst T. and Prokop, Mathias and Jacobs, Colin'
$ws.Range("I121").Value2 = 2023

$ws.Range("A122").Value2 = 1746
$ws.Range("B122").Value2 = 'Graa23a'
$ws.Range("C122").Value2 = 'article'
$ws.Range("D122").Value2 = 'Lumbar spine segmentation in MR images: a dataset and a public benchmark'
$ws.Range("E122").Value2 = 'van der Graaf, Jasper W. and van Hooff, Miranda L. and Buckens, Constantinus F. M. and Rutten, Matthieu and van Susante, Job L. C. and Kroeze, Robert Jan and de Kleuver, Marinus and van Ginneken, Bram and Lessmann, Nikolas'
$ws.Range("I122").Value2 = 2023

$ws.Range("A123").Value2 = 1747
$ws.Range("B123").Value2 = 'Thij23'
$ws.Range("C123").Value2 = 'article'
$ws.Range("D123").Value2 = 'Radiomics based automated quality assessment for T2W prostate MR images'
$ws.Range("E123").Value2 = 'Thijssen, Linda C.P. and de Rooij, Maarten and Barentsz, Jelle O. and Huisman, Henkjan J.'
$ws.Range("I123").Value2 = 2023
